# Added feature to extract .msg attachments - populate form-submission rows
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that get a simple "Yes" marker for each new submission row
$yesCols = @("N", "R", "W", "Y", "AQ", "BN", "BU", "BW")

for ($row = 2; $row -le 5; $row++) {
    foreach ($col in $yesCols) {
        $ws.Range("$col$row").Value = "Yes"
    }
    # Force the date-like string to be stored as literal text (not an
    # auto-converted date serial): apply a text format just long enough to
    # enter the value, then drop back to the default "Normal" style so the
    # cell is left without any explicit formatting, matching the other
    # plain data cells in the row.
    $ws.Range("AI$row").NumberFormat = "@"
    $ws.Range("AI$row").Value = "2022/01/08"
    $ws.Range("AI$row").Style = "Normal"
    $ws.Range("AJ$row").Value = "Chris Gryzen"
    $ws.Range("AK$row").Value = "Josh Gryzen"
    $ws.Range("AL$row").Value = "Gryzen"
    $ws.Range("CD$row").Value = 17
}
